$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2500913
$ws.Range("J17").Value = 2500913
$ws.Range("L17").Value = 7502739
$ws.Range("N17").Value = -7503075
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
# Row 137
$ws.Range("H137").Value = 6048.1
$ws.Range("J137").Value = 8096.1665
$ws.Range("L137").Value = 24288.4995
$ws.Range("N137").Value = -29388.4995
# Row 138
$ws.Range("H138").Value = 3155.457
$ws.Range("I138").Value = 2455.75
$ws.Range("J138").Value = 3300.224
$ws.Range("K138").Value = 7367.25
$ws.Range("L138").Value = 9900.672
$ws.Range("M138").Value = -2227.25
$ws.Range("N138").Value = -20180.672

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 55557376
$ws.Range("I61").Value = 71430330
$ws.Range("K61").Value = 71430330
$ws.Range("M61").Value = -71430118
# Row 105
$ws.Range("H105").Value = 6669.3335
$ws.Range("I105").Value = 6669.3335
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 6669.3335
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3175.3335
$ws.Range("N105").ClearContents()
# Row 136
$ws.Range("H136").Value = 55557376
$ws.Range("I136").Value = 71430330
$ws.Range("K136").Value = 214290990
$ws.Range("M136").Value = -214288440

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 17382.285
$ws.Range("I105").Value = 28280
$ws.Range("J105").Value = 2852
$ws.Range("K105").Value = 28280
$ws.Range("L105").Value = 2852
$ws.Range("M105").Value = -26533
$ws.Range("N105").Value = -6346
# Row 107
$ws.Range("H107").Value = 1523
$ws.Range("I107").Value = 1112.0952
$ws.Range("J107").Value = 2307.4546
$ws.Range("K107").Value = 1112.0952
$ws.Range("L107").Value = 2307.4546
$ws.Range("M107").Value = 807.9048
$ws.Range("N107").Value = -6147.4546
# Row 110
$ws.Range("H110").Value = 49000
$ws.Range("J110").Value = 49000
$ws.Range("L110").Value = 49000
$ws.Range("N110").Value = -57180

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 94555000
$ws.Range("I4").Value = 13333340
$ws.Range("J4").Value = 125013130
$ws.Range("K4").Value = 13333340
$ws.Range("L4").Value = 125013130
$ws.Range("M4").Value = -13333228
$ws.Range("N4").Value = -125013354
# Row 31
$ws.Range("H31").Value = 4903966
$ws.Range("I31").Value = 1882.0526
$ws.Range("J31").Value = 7814578.5
$ws.Range("K31").Value = 1882.0526
$ws.Range("L31").Value = 7814578.5
$ws.Range("M31").Value = -1587.0526
$ws.Range("N31").Value = -7815168.5
# Row 34
$ws.Range("H34").Value = 4903966
$ws.Range("I34").Value = 1882.0526
$ws.Range("J34").Value = 7814578.5
$ws.Range("K34").Value = 1882.0526
$ws.Range("L34").Value = 7814578.5
$ws.Range("M34").Value = -1680.0526
$ws.Range("N34").Value = -7814982.5
# Row 41
$ws.Range("H41").Value = 24873.334
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 43972
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 43972
$ws.Range("M41").Value = -572
$ws.Range("N41").Value = -44828
# Row 51
$ws.Range("H51").Value = 29753.545
$ws.Range("J51").Value = 29753.545
$ws.Range("L51").Value = 29753.545
$ws.Range("N51").Value = -31225.545
# Row 61
$ws.Range("H61").Value = 29753.545
$ws.Range("J61").Value = 29753.545
$ws.Range("L61").Value = 29753.545
$ws.Range("N61").Value = -30449.545
# Row 68
$ws.Range("H68").Value = 99500
$ws.Range("J68").Value = 99500
$ws.Range("L68").Value = 99500
$ws.Range("N68").Value = -100998
# Row 71
$ws.Range("H71").Value = 99500
$ws.Range("J71").Value = 99500
$ws.Range("L71").Value = 298500
$ws.Range("N71").Value = -305988
# Row 74
$ws.Range("H74").Value = 39555
$ws.Range("J74").Value = 39555
$ws.Range("L74").Value = 39555
$ws.Range("N74").Value = -41303
# Row 77
$ws.Range("H77").Value = 39555
$ws.Range("J77").Value = 39555
$ws.Range("L77").Value = 118665
$ws.Range("N77").Value = -127401
# Row 132
$ws.Range("H132").Value = 43651.418
$ws.Range("I132").Value = 47392.203
$ws.Range("J132").Value = 2502.75
$ws.Range("K132").Value = 142176.609
$ws.Range("L132").Value = 7508.25
$ws.Range("M132").Value = -139646.609
$ws.Range("N132").Value = -12568.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 38036330
$ws.Range("I4").Value = 53495020
$ws.Range("J4").Value = 20001188
$ws.Range("K4").Value = 160485060
$ws.Range("L4").Value = 60003564
$ws.Range("M4").Value = -160484948
$ws.Range("N4").Value = -60003788
# Row 107
$ws.Range("H107").Value = 1193.409
$ws.Range("J107").Value = 1956.5
$ws.Range("L107").Value = 5869.5
$ws.Range("N107").Value = -9709.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 1515287.2
$ws.Range("I2").Value = 2272798.8
$ws.Range("J2").Value = 264.36365
$ws.Range("K2").Value = 2272798.8
$ws.Range("L2").Value = 264.36365
$ws.Range("M2").Value = -2272685.8
$ws.Range("N2").Value = -490.36365
# Row 126
$ws.Range("H126").Value = 95182.82
$ws.Range("I126").Value = 203402.2
$ws.Range("K126").Value = 610206.6000000001
$ws.Range("M126").Value = -607736.6000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 7917.643
$ws.Range("J40").Value = 7559
$ws.Range("L40").Value = 7559
$ws.Range("N40").Value = -7831
# Row 61
$ws.Range("H61").Value = 2520
$ws.Range("I61").Value = 1594.7333
$ws.Range("K61").Value = 1594.7333
$ws.Range("M61").Value = -1392.7333
# Row 63
$ws.Range("H63").Value = 66925.164
$ws.Range("J63").Value = 61332
$ws.Range("L63").Value = 61332
$ws.Range("N63").Value = -62830
# Row 66
$ws.Range("H66").Value = 66925.164
$ws.Range("J66").Value = 61332
$ws.Range("L66").Value = 183996
$ws.Range("N66").Value = -191484
# Row 113
$ws.Range("H113").Value = 2520
$ws.Range("I113").Value = 1594.7333
$ws.Range("K113").Value = 1594.7333
$ws.Range("M113").Value = 575.2666999999999
# Row 131
$ws.Range("H131").Value = 73463.664
$ws.Range("J131").Value = 73463.664
$ws.Range("L131").Value = 73463.664
$ws.Range("N131").Value = -83543.664
# Row 136
$ws.Range("H136").Value = 1253924.5
$ws.Range("I136").Value = 1820981.2
$ws.Range("J136").Value = 6399.8
$ws.Range("K136").Value = 5462943.6
$ws.Range("L136").Value = 19199.4
$ws.Range("M136").Value = -5460393.6
$ws.Range("N136").Value = -24299.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 96769.81
$ws.Range("I2").Value = 133146.89
$ws.Range("K2").Value = 133146.89
$ws.Range("M2").Value = -133034.89
# Row 75
$ws.Range("H75").Value = 49997
$ws.Range("J75").Value = 49999.332
$ws.Range("L75").Value = 49999.332
$ws.Range("N75").Value = -51871.332
# Row 78
$ws.Range("H78").Value = 49997
$ws.Range("J78").Value = 49999.332
$ws.Range("L78").Value = 149997.996
$ws.Range("N78").Value = -159357.996
# Row 100
$ws.Range("H100").Value = 56111816
$ws.Range("I100").Value = 67333880
$ws.Range("J100").Value = 1467.6666
$ws.Range("K100").Value = 134667760
$ws.Range("L100").Value = 2935.3332
$ws.Range("M100").Value = -134667219
$ws.Range("N100").Value = -4017.3332
# Row 104
$ws.Range("H104").Value = 37013.6
$ws.Range("J104").Value = 37013.6
$ws.Range("L104").Value = 37013.6
$ws.Range("N104").Value = -44001.6
# Row 136
$ws.Range("H136").Value = 2521.0833
$ws.Range("I136").Value = 3183.8333
$ws.Range("J136").Value = 1858.3334
$ws.Range("K136").Value = 9551.499899999999
$ws.Range("L136").Value = 5575.0002
$ws.Range("M136").Value = -7001.499899999999
$ws.Range("N136").Value = -10675.0002
# Row 138
$ws.Range("H138").Value = 130429
$ws.Range("J138").Value = 130429
$ws.Range("L138").Value = 130429
$ws.Range("N138").Value = -140709
